$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format before assigning so that numeric-looking
    # strings (e.g. "1.00", "573.84") are preserved exactly as text instead
    # of being auto-converted into a Number by Excel. ClearFormats afterwards
    # removes the temporary "@" number-format style so the cell keeps the
    # same (unstyled) appearance as before.
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "67.180.63"
$ws.Range("E2").Value2 = "  -0.19%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.106.38"
$ws.Range("E3").Value2 = "  -0.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  -0.14%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "573.84"
$ws.Range("E5").Value2 = "  -1.06%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "177.97"
$ws.Range("E6").Value2 = "  +2.97%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value2 = "  -0.08%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.105.73"
$ws.Range("E8").Value2 = "  +0.12%  "

# Row 9 - XRP
$ws.Range("E9").Value2 = "  -1.27%  "

# Row 10 - Toncoin
$ws.Range("E10").Value2 = "  -1.75%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value2 = "  -0.57%  "

# Row 12 - Cardano
$ws.Range("E12").Value2 = "  -1.72%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value2 = "  -2.53%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "36.07"
$ws.Range("E14").Value2 = "  -1.56%  "

# Row 15 - TRON
$ws.Range("E15").Value2 = "  +0.01%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D16") "3.624.36"
$ws.Range("E16").Value2 = "  -0.06%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "67.079.98"
$ws.Range("E17").Value2 = "  -0.36%  "

# Row 19 - WrappedEther
Set-TextValue $ws.Range("D19") "3.107.64"
$ws.Range("E19").Value2 = "  -0.14%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "16.69"
$ws.Range("E20").Value2 = "  +0.18%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "493.09"
$ws.Range("E21").Value2 = "  +0.30%  "

# Row 22 - Uniswap
$ws.Range("E22").Value2 = "  -0.09%  "

# Row 23 - Polygon
$ws.Range("E23").Value2 = "  -1.61%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "83.62"
$ws.Range("E24").Value2 = "  -0.35%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D25") "12.60"
$ws.Range("E25").Value2 = "  -3.52%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value2 = "  -0.78%  "

# Row 27 - RenderToken
$ws.Range("E27").Value2 = "  -4.30%  "

# Row 29 - NEARProtocol
Set-TextValue $ws.Range("D29") "7.93"
$ws.Range("E29").Value2 = "  +1.01%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "2.31"
$ws.Range("E30").Value2 = "  -1.21%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value2 = "  -2.61%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "28.11"
$ws.Range("E32").Value2 = "  -0.68%  "

# Row 33 - Hedera
$ws.Range("E33").Value2 = "  -1.52%  "

# Row 34 - PEPE
Set-TextValue $ws.Range("D34") "0.0₃0940"
$ws.Range("E34").Value2 = "  -0.36%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value2 = "  -0.07%  "

# Row 36 - Arweave
Set-TextValue $ws.Range("D36") "47.63"
$ws.Range("E36").Value2 = "  +2.12%  "

# Row 37 / 38 - Mantle and Filecoin swap ranking order
$ws.Range("B37").Value2 = "Filecoin"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D37") "5.58"
$ws.Range("E37").Value2 = "  -3.64%  "

$ws.Range("B38").Value2 = "Mantle"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D38") "0.945"
$ws.Range("E38").Value2 = "  -2.72%  "

# Row 39 - TheGraph
Set-TextValue $ws.Range("D39") "0.312"
$ws.Range("E39").Value2 = "  +1.71%  "

# Row 40 - OKB
Set-TextValue $ws.Range("D40") "49.11"
$ws.Range("E40").Value2 = "  -1.50%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("D41") "2.01"
$ws.Range("E41").Value2 = "  -0.85%  "

# Row 42 - Kaspa
$ws.Range("E42").Value2 = "  -0.45%  "

# Row 43 - Cosmos
$ws.Range("E43").Value2 = "  -1.92%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value2 = "  +6.48%  "

# Row 45 - Maker
Set-TextValue $ws.Range("D45") "2.795.09"
$ws.Range("E45").Value2 = "  -0.12%  "

# Row 46 - Bittensor
Set-TextValue $ws.Range("D46") "370.49"
$ws.Range("E46").Value2 = "  -3.49%  "

# Row 47 - VeChain
$ws.Range("E47").Value2 = "  -1.67%  "

# Row 48 - Monero
Set-TextValue $ws.Range("D48") "135.78"
$ws.Range("E48").Value2 = "  +0.34%  "

# Row 49 - USDe
$ws.Range("E49").Value2 = "  +0.01%  "

# Row 50 - InjectiveProtocol
Set-TextValue $ws.Range("D50") "25.44"
$ws.Range("E50").Value2 = "  +1.93%  "

# Row 51 - ThetaToken
Set-TextValue $ws.Range("D51") "2.28"
$ws.Range("E51").Value2 = "  +4.00%  "
